# Apply "Add data for 2022-06-29" update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet to reflect new "through" date
$ws.Name = "Through 2022-06-21"

# Update the June row label (shared string) to reflect new "through" date
$ws.Range("A7").Value = "June (through 06-21)"

# Update June (row 7) figures for columns D, E, G, H, I (2017, 2018, 2020, 2021, 2022)
$ws.Range("D7").Value = 51
$ws.Range("E7").Value = 39
$ws.Range("G7").Value = 82
$ws.Range("H7").Value = 82
$ws.Range("I7").Value = 101

# Update Total row (row 8) figures for the same columns
$ws.Range("D8").Value = 367
$ws.Range("E8").Value = 334
$ws.Range("G8").Value = 440
$ws.Range("H8").Value = 713
$ws.Range("I8").Value = 764
